$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 5: Best Up College Information Workshop ---
$ws.Range("A5").Value = "BestUpInformationWorkshop"
$ws.Range("B5").Value = "Best Up College Information Workshop"
$ws.Range("C5").Value = 2013
$ws.Range("D5").Value = "Information"
$ws.Range("E5").Value = 19
$ws.Range("F5").Value = "A group of researchers  studied the effect of providing information about the benefits of college education to high school students one year prior to graduation."
$ws.Range("F5").WrapText = $true

# --- Add new row 6: Mentoring Program Balu und Du ---
$ws.Range("A6").Value = "mentoringBalu"
$ws.Range("B6").Value = "Mentoring Program Balu und Du"
$ws.Range("C6").Value = 2011
$ws.Range("D6").Value = "Mentoring"
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = "Elementary School students were assigned mentors who should encourage the acquisition of new ideas and skills."
$ws.Range("F6").WrapText = $true

# --- Updated row heights (graph/layout refresh) ---
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(6).RowHeight = 45

# --- Selection ends on the newly added row ---
[void]$ws.Range("A6").Select()
